$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CP_CONCLNEG_001): "Resultado Obtenido" text now references the
# Contenido clases de negocio view instead of Autodiagnostico.
$ws.Range("I2").Value = "La vista Contenido clases de negocio se cargó sin errores"

# Row 3 (CP_CONCLNEG_002): fill in the previously empty "Resultado Esperado"
# / "Resultado Obtenido" columns.
$ws.Range("H3").Value = "El sistema permite seleccionar la entidad y mostrar correctamente la vista (Modelos)"
$ws.Range("I3").Value = "la vista (modelos) se visualiza correctamente"

# Row 4 (CP_CONCLNEG_003): fill in the previously empty "Resultado Esperado".
$ws.Range("H4").Value = "El sistema la creacion de un modelo correctamente"

# Row 4 grew taller to fit its new text.
$ws.Rows(4).RowHeight = 216.75

# Update the view: drop the old scroll position / selection and select I4.
$ws.Range("I4").Select()
